$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version bump 0.1.0 -> 0.1.1
$meta.Cells.Item(3, 2).Value = "0.1.1"

# Date bump
$meta.Cells.Item(8, 2).Value = "2023-06-02T12:02:38+02:00"

# Context now applies to five resource types instead of one. Row 20
# already holds the (only) "Context" row (element:PractitionerRole);
# clone its formatting into four new rows below it, then fill in the
# Property/Value pairs for all five contexts (Location, Organization,
# Practitioner, PractitionerRole, HealthcareService).
$meta.Range("A19:B19").Copy()
$meta.Range("A21:B24").PasteSpecial(-4122)

$meta.Cells.Item(20, 1).Value = "Context"
$meta.Cells.Item(20, 2).Value = "element:Location"

$meta.Cells.Item(21, 1).Value = "Context"
$meta.Cells.Item(21, 2).Value = "element:Organization"

$meta.Cells.Item(22, 1).Value = "Context"
$meta.Cells.Item(22, 2).Value = "element:Practitioner"

$meta.Cells.Item(23, 1).Value = "Context"
$meta.Cells.Item(23, 2).Value = "element:PractitionerRole"

$meta.Cells.Item(24, 1).Value = "Context"
$meta.Cells.Item(24, 2).Value = "element:HealthcareService"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type(s): date -> dateTime
$elements.Cells.Item(6, 11).Value = "dateTime`n"
